# Apply the latest crypto price/volume snapshot scraped by the
# GitHub Actions job (cryptos.xlsx refresh). Price/volume cells are
# stored as plain text (e.g. "30.666.25", "  +2.41%  "), so every
# write below keeps that text formatting intact rather than letting
# Excel reinterpret the price strings as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.666.25"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.74"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.67"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4925"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2948"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06793"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.886.76"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.14"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07244"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.15"
$ws.Range("E13").Value = "  +5.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.057"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6763"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.650.35"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007963"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.132.34"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.821"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "183.89"
$ws.Range("E23").Value = "  +29.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.034"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.329"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.85"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.901"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.400"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.284"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08986"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.989"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05187"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7396"
$ws.Range("E34").Value = "  +3.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.742"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01835"
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.668"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.144"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9369"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4408"
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.19"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.760"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.584"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1331"
$ws.Range("E46").Value = "  +4.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05844"
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.432"
$ws.Range("E48").Value = "  +7.17%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.606"
$ws.Range("E49").Value = "  +3.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3922"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.42"
$ws.Range("E51").Value = "  +2.37%  "
